# The document contains a paragraph with a Word field (MERGEFIELD-like
# construct built from fldChar/instrText runs: " m: ") immediately followed
# by a bold red run of error text. The commit replaces the field-code
# construct with plain literal text runs "{", "m", ":}" (mirroring the
# field's instrText content, now rendered as ordinary template-syntax
# text instead of a real Word field) and prefixes the error message with
# "    <---".

$d = $word.ActiveDocument

# Locate the field (there is exactly one in this document) and remember
# where it starts so we can insert the replacement text in its place.
$f = $d.Fields.Item(1)
$fieldStart = $f.Code.Start - 1

# Remove the field (fldChar begin/end + instrText runs) entirely.
$f.Delete()

# Insert three plain-text runs "{", "m", ":}" exactly where the field used
# to be - split into separate runs/w:t elements to mirror the original
# instrText segmentation ("m" and ": ").
$insertRange = $d.Range($fieldStart, $fieldStart)
$insertRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t xml:space="preserve">:}</w:t></w:r></w:p>')

# Prefix the (still bold/red) error message run with four spaces and "<---".
$d.Content.Find.Execute(
    "Invalid query statement: Expression", $true, $false, $false, $false,
    $false, $true, 1, $false, "    <---Invalid query statement: Expression",
    2)
